# Refresh cryptos list values (price & 1h volume change) plus a few reordered rows,
# matching the latest scrape from coinranking.com. Text values are written with a
# leading apostrophe so Excel keeps storing them as text (matching the original
# inlineStr cells) instead of auto-converting number-looking strings like "1.00".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.001.71"
$ws.Range("E2").Value = "'  -5.35%  "
$ws.Range("D3").Value = "'3.694.19"
$ws.Range("E3").Value = "'  -5.11%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  -0.17%  "
$ws.Range("D5").Value = "'583.71"
$ws.Range("E5").Value = "'  -2.42%  "
$ws.Range("D6").Value = "'179.33"
$ws.Range("E6").Value = "'  +6.93%  "
$ws.Range("D7").Value = "'3.689.07"
$ws.Range("E7").Value = "'  -5.18%  "
$ws.Range("D8").Value = "'0.627"
$ws.Range("E8").Value = "'  -6.44%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "'  -0.13%  "
$ws.Range("D10").Value = "'0.709"
$ws.Range("E10").Value = "'  -7.30%  "
$ws.Range("D11").Value = "'0.161"
$ws.Range("E11").Value = "'  -9.63%  "
$ws.Range("D12").Value = "'53.83"
$ws.Range("E12").Value = "'  -0.44%  "
$ws.Range("D13").Value = "'0.0000290"
$ws.Range("E13").Value = "'  -10.34%  "
$ws.Range("D14").Value = "'10.37"
$ws.Range("E14").Value = "'  -8.57%  "
$ws.Range("D15").Value = "'4.352.88"
$ws.Range("E15").Value = "'  -4.06%  "
$ws.Range("D16").Value = "'3.704.57"
$ws.Range("E16").Value = "'  -5.53%  "
$ws.Range("D17").Value = "'19.37"
$ws.Range("E17").Value = "'  -8.12%  "
$ws.Range("E18").Value = "'  -2.88%  "
$ws.Range("D19").Value = "'12.78"
$ws.Range("E19").Value = "'  -8.35%  "
$ws.Range("E20").Value = "'  -7.86%  "
$ws.Range("D21").Value = "'67.864.87"
$ws.Range("E21").Value = "'  -5.71%  "
$ws.Range("D22").Value = "'406.59"
$ws.Range("E22").Value = "'  -6.82%  "
$ws.Range("D23").Value = "'4.49"
$ws.Range("E23").Value = "'  -5.58%  "
$ws.Range("D24").Value = "'88.14"
$ws.Range("E24").Value = "'  -6.66%  "
$ws.Range("D25").Value = "'3.02"
$ws.Range("E25").Value = "'  -8.57%  "
$ws.Range("D26").Value = "'12.73"
$ws.Range("E26").Value = "'  -8.14%  "
$ws.Range("D27").Value = "'10.96"
$ws.Range("E27").Value = "'  -0.42%  "
$ws.Range("D28").Value = "'3.86"
$ws.Range("E28").Value = "'  -7.12%  "
$ws.Range("D29").Value = "'6.07"
$ws.Range("E29").Value = "'  +2.14%  "
$ws.Range("D30").Value = "'9.46"
$ws.Range("E30").Value = "'  -7.40%  "
$ws.Range("D31").Value = "'32.39"
$ws.Range("E31").Value = "'  -8.00%  "
$ws.Range("D32").Value = "'7.48"
$ws.Range("E32").Value = "'  -6.58%  "
$ws.Range("D33").Value = "'12.44"
$ws.Range("E33").Value = "'  -8.84%  "
$ws.Range("B34").Value = "'OKB"
$ws.Range("C34").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'65.32"
$ws.Range("E34").Value = "'  -4.63%  "
$ws.Range("B35").Value = "'Hedera"
$ws.Range("C35").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.116"
$ws.Range("E35").Value = "'  -7.85%  "
$ws.Range("B36").Value = "'Bittensor"
$ws.Range("C36").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").Value = "'598.40"
$ws.Range("E36").Value = "'  -4.02%  "
$ws.Range("B37").Value = "'InjectiveProtocol"
$ws.Range("C37").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "'42.86"
$ws.Range("E37").Value = "'  -18.30%  "
$ws.Range("D38").Value = "'0.0₃0890"
$ws.Range("E38").Value = "'  -9.20%  "
$ws.Range("E39").Value = "'  +0.06%  "
$ws.Range("D40").Value = "'0.396"
$ws.Range("E40").Value = "'  -5.79%  "
$ws.Range("E41").Value = "'  -0.36%  "
$ws.Range("E42").Value = "'  -4.91%  "
$ws.Range("D43").Value = "'2.76"
$ws.Range("E43").Value = "'  +4.38%  "
$ws.Range("D44").Value = "'2.98"
$ws.Range("E44").Value = "'  -10.14%  "
$ws.Range("D45").Value = "'2.94"
$ws.Range("E45").Value = "'  -8.43%  "
$ws.Range("D46").Value = "'0.0433"
$ws.Range("E46").Value = "'  -8.19%  "
$ws.Range("D47").Value = "'9.23"
$ws.Range("E47").Value = "'  -10.49%  "
$ws.Range("D48").Value = "'2.797.40"
$ws.Range("E48").Value = "'  -2.94%  "
$ws.Range("B49").Value = "'WEMIXToken"
$ws.Range("C49").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'2.70"
$ws.Range("E49").Value = "'  -5.45%  "
$ws.Range("B50").Value = "'Stellar"
$ws.Range("C50").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.133"
$ws.Range("E50").Value = "'  -7.80%  "
$ws.Range("D51").Value = "'3.08"
$ws.Range("E51").Value = "'  -7.68%  "
